$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value2 = 67114.17999999999
$ws.Range("I70").Value2 = 202040.2
$ws.Range("J70").Value2 = 10895
$ws.Range("K70").Value2 = 606120.6000000001
$ws.Range("L70").Value2 = 32685
$ws.Range("M70").Value2 = -605850.6000000001
$ws.Range("N70").Value2 = -33225
# Row 73
$ws.Range("H73").Value2 = 67114.17999999999
$ws.Range("I73").Value2 = 202040.2
$ws.Range("J73").Value2 = 10895
$ws.Range("K73").Value2 = 606120.6000000001
$ws.Range("L73").Value2 = 32685
$ws.Range("M73").Value2 = -605184.6000000001
$ws.Range("N73").Value2 = -34557
# Row 86
$ws.Range("H86").Value2 = 3773.4285
$ws.Range("I86").Value2 = 2078.8333
$ws.Range("K86").Value2 = 2078.8333
$ws.Range("M86").Value2 = -955.8332999999998
# Row 89
$ws.Range("H89").Value2 = 3773.4285
$ws.Range("I89").Value2 = 2078.8333
$ws.Range("K89").Value2 = 10394.1665
$ws.Range("M89").Value2 = -4778.166499999999
# Row 92
$ws.Range("H92").Value2 = 1003
$ws.Range("I92").Value2 = 901.0526
$ws.Range("K92").Value2 = 901.0526
$ws.Range("M92").Value2 = 346.9474
# Row 113
$ws.Range("H113").Value2 = 4946.5
$ws.Range("J113").Value2 = 5413.7856
$ws.Range("L113").Value2 = 5413.7856
$ws.Range("N113").Value2 = -11921.7856
# Row 116
$ws.Range("H116").Value2 = 2592.0833
$ws.Range("I116").Value2 = 2583.3333
$ws.Range("K116").Value2 = 2583.3333
$ws.Range("M116").Value2 = 858.6667000000002
# Row 132
$ws.Range("H132").Value2 = 7583.6045
$ws.Range("I132").Value2 = 3945
$ws.Range("J132").Value2 = 26296.428
$ws.Range("K132").Value2 = 11835
$ws.Range("L132").Value2 = 78889.284
$ws.Range("M132").Value2 = -9305
$ws.Range("N132").Value2 = -83949.284
# Row 138
$ws.Range("H138").Value2 = 1102834.8
$ws.Range("I138").Value2 = 334782.34
$ws.Range("J138").Value2 = 1432000
$ws.Range("K138").Value2 = 1004347.02
$ws.Range("L138").Value2 = 4296000
$ws.Range("M138").Value2 = -999207.02
$ws.Range("N138").Value2 = -4306280

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value2 = 50000
$ws.Range("I5").Value2 = 50000
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 50000
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = -49888
$ws.Range("N5").Value2 = $null
# Row 11
$ws.Range("H11").Value2 = 14537.6
$ws.Range("I11").Value2 = 9566.333000000001
$ws.Range("J11").Value2 = 21994.5
$ws.Range("K11").Value2 = 9566.333000000001
$ws.Range("L11").Value2 = 21994.5
$ws.Range("M11").Value2 = -9422.333000000001
$ws.Range("N11").Value2 = -22282.5
# Row 45
$ws.Range("H45").Value2 = 2023.52
$ws.Range("I45").Value2 = 1851.6522
$ws.Range("K45").Value2 = 1851.6522
$ws.Range("M45").Value2 = -1474.6522
# Row 111
$ws.Range("H111").Value2 = 87665
$ws.Range("J111").Value2 = 87665
$ws.Range("L111").Value2 = 87665
$ws.Range("N111").Value2 = -95845
# Row 132
$ws.Range("H132").Value2 = 5577.9
$ws.Range("I132").Value2 = 1730.7333
$ws.Range("J132").Value2 = 17119.4
$ws.Range("K132").Value2 = 5192.199900000001
$ws.Range("L132").Value2 = 51358.2
$ws.Range("M132").Value2 = -2662.199900000001
$ws.Range("N132").Value2 = -56418.2

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value2 = 50000
$ws.Range("I4").Value2 = 50000
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 50000
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = -49885
$ws.Range("N4").Value2 = $null
# Row 20
$ws.Range("H20").Value2 = 2263.182
$ws.Range("I20").Value2 = 2108.4614
$ws.Range("J20").Value2 = 2486.6667
$ws.Range("K20").Value2 = 2108.4614
$ws.Range("L20").Value2 = 2486.6667
$ws.Range("M20").Value2 = -1861.4614
$ws.Range("N20").Value2 = -2980.6667
# Row 134
$ws.Range("H134").Value2 = 6063.2583
$ws.Range("I134").Value2 = 2223.1428
$ws.Range("K134").Value2 = 6669.428400000001
$ws.Range("M134").Value2 = -4134.428400000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value2 = 622.125
$ws.Range("I22").Value2 = 586.7273
$ws.Range("J22").Value2 = 700
$ws.Range("K22").Value2 = 586.7273
$ws.Range("L22").Value2 = 700
$ws.Range("M22").Value2 = -236.7273
$ws.Range("N22").Value2 = -1400
# Row 31
$ws.Range("H31").Value2 = 2281.2295
$ws.Range("I31").Value2 = 1731.14
$ws.Range("J31").Value2 = 4781.636
$ws.Range("K31").Value2 = 1731.14
$ws.Range("L31").Value2 = 4781.636
$ws.Range("M31").Value2 = -1436.14
$ws.Range("N31").Value2 = -5371.636
# Row 34
$ws.Range("H34").Value2 = 2281.2295
$ws.Range("I34").Value2 = 1731.14
$ws.Range("J34").Value2 = 4781.636
$ws.Range("K34").Value2 = 1731.14
$ws.Range("L34").Value2 = 4781.636
$ws.Range("M34").Value2 = -1529.14
$ws.Range("N34").Value2 = -5185.636
# Row 42
$ws.Range("H42").Value2 = 4000
$ws.Range("J42").Value2 = 4000
$ws.Range("L42").Value2 = 4000
$ws.Range("N42").Value2 = -5186
# Row 58
$ws.Range("H58").Value2 = 2103.818
$ws.Range("I58").Value2 = 2113.3572
$ws.Range("K58").Value2 = 2113.3572
$ws.Range("M58").Value2 = -1910.3572
# Row 97
$ws.Range("H97").Value2 = 89197
$ws.Range("J97").Value2 = 89197
$ws.Range("L97").Value2 = 89197
$ws.Range("N97").Value2 = -91179
# Row 122
$ws.Range("H122").Value2 = 3216.3948
$ws.Range("I122").Value2 = 3174.138
$ws.Range("K122").Value2 = 9522.414000000001
$ws.Range("M122").Value2 = -7072.414000000001
# Row 132
$ws.Range("H132").Value2 = 932582.6
$ws.Range("I132").Value2 = 1178887.1
$ws.Range("K132").Value2 = 3536661.3
$ws.Range("M132").Value2 = -3534131.3
# Row 134
$ws.Range("H134").Value2 = 2455.7925
$ws.Range("I134").Value2 = 1692.2046
$ws.Range("J134").Value2 = 6188.8887
$ws.Range("K134").Value2 = 5076.6138
$ws.Range("L134").Value2 = 18566.6661
$ws.Range("M134").Value2 = -2541.6138
$ws.Range("N134").Value2 = -23636.6661
# Row 136
$ws.Range("H136").Value2 = 2103.818
$ws.Range("I136").Value2 = 2113.3572
$ws.Range("K136").Value2 = 6340.071599999999
$ws.Range("M136").Value2 = -3790.071599999999
# Row 141
$ws.Range("H141").Value2 = 124975.78
$ws.Range("J141").Value2 = 124975.78
$ws.Range("L141").Value2 = 124975.78
$ws.Range("N141").Value2 = -135335.78

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value2 = 245.66667
$ws.Range("J92").Value2 = 246.25
$ws.Range("L92").Value2 = 738.75
$ws.Range("N92").Value2 = -3234.75
# Row 97
$ws.Range("H97").Value2 = 1405.6
$ws.Range("I97").Value2 = 1622.2858
$ws.Range("K97").Value2 = 4866.857400000001
$ws.Range("M97").Value2 = -4370.857400000001
# Row 113
$ws.Range("H113").Value2 = 990.2857
$ws.Range("J113").Value2 = 1057.7894
$ws.Range("L113").Value2 = 3173.3682
$ws.Range("N113").Value2 = -7513.3682
# Row 131
$ws.Range("H131").Value2 = 1881.3125
$ws.Range("J131").Value2 = 1881.3125
$ws.Range("L131").Value2 = 5643.9375
$ws.Range("N131").Value2 = -15723.9375
# Row 133
$ws.Range("H133").Value2 = 5717.8
$ws.Range("I133").Value2 = 5022.25
$ws.Range("K133").Value2 = 15066.75
$ws.Range("M133").Value2 = -10006.75
# Row 140
$ws.Range("H140").Value2 = 2885.724
$ws.Range("J140").Value2 = 5312.5
$ws.Range("L140").Value2 = 15937.5
$ws.Range("N140").Value2 = -26297.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value2 = 4042.7144
$ws.Range("I70").Value2 = 2824.75
$ws.Range("J70").Value2 = 5666.6665
$ws.Range("K70").Value2 = 2824.75
$ws.Range("L70").Value2 = 5666.6665
$ws.Range("M70").Value2 = -2554.75
$ws.Range("N70").Value2 = -6206.6665
# Row 73
$ws.Range("H73").Value2 = 4042.7144
$ws.Range("I73").Value2 = 2824.75
$ws.Range("J73").Value2 = 5666.6665
$ws.Range("K73").Value2 = 2824.75
$ws.Range("L73").Value2 = 5666.6665
$ws.Range("M73").Value2 = -1888.75
$ws.Range("N73").Value2 = -7538.6665
# Row 102
$ws.Range("H102").Value2 = 65644.25
$ws.Range("I102").Value2 = 4130
$ws.Range("K102").Value2 = 4130
$ws.Range("M102").Value2 = -2508
# Row 132
$ws.Range("H132").Value2 = 3465.3823
$ws.Range("I132").Value2 = 3472.6875
$ws.Range("K132").Value2 = 10418.0625
$ws.Range("M132").Value2 = -7888.0625

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value2 = 6640.913
$ws.Range("I46").Value2 = 3087.5
$ws.Range("J46").Value2 = 8536.066000000001
$ws.Range("K46").Value2 = 3087.5
$ws.Range("L46").Value2 = 8536.066000000001
$ws.Range("M46").Value2 = -2899.5
$ws.Range("N46").Value2 = -8912.066000000001
# Row 122
$ws.Range("H122").Value2 = 3588.96
$ws.Range("I122").Value2 = 2976.6667
$ws.Range("K122").Value2 = 8930.000100000001
$ws.Range("M122").Value2 = -6480.000100000001
# Row 132
$ws.Range("H132").Value2 = 3160.9736
$ws.Range("J132").Value2 = 3472
$ws.Range("L132").Value2 = 10416
$ws.Range("N132").Value2 = -15476

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value2 = 92944
$ws.Range("J46").Value2 = 92944
$ws.Range("L46").Value2 = 92944
$ws.Range("N46").Value2 = -93406
# Row 132
$ws.Range("H132").Value2 = 2628.7917
$ws.Range("I132").Value2 = 2663.2273
$ws.Range("K132").Value2 = 7989.6819
$ws.Range("M132").Value2 = -5459.6819
# Row 134
$ws.Range("H134").Value2 = 92944
$ws.Range("J134").Value2 = 92944
$ws.Range("L134").Value2 = 278832
$ws.Range("N134").Value2 = -283902
